# Qvx Reader update in progress
# In the process of reading Date/Time/Timestamp/Interval values
#
# On the "Tasks" sheet, the Reader/Writer columns for the DATE, TIME,
# TIMESTAMP and INTERVAL rows (rows 8-11) move from "not started" (red)
# to "in progress" (yellow): the Reader column now notes that a reg-exp
# function is used, and the Writer column is filled in with the actual
# format strings used for each data type. The two old placeholder notes
# ("Unsure of how FieldAttributes.fmt works" / "Unsure of what this means
# ...") are replaced by real content; a new "Unsure of how to use Excel to
# generate this data type" note is added for the TIME row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Reader column: all four rows now read "Reg-exp function"
$ws.Range("B8:B11").Value = "Reg-exp function"

# Writer column: the concrete format string used for each type
$ws.Range("C8").Value  = "M/D/YYYY"
$ws.Range("C9").Value  = "Unsure of how to use Excel to generate this data type"
$ws.Range("C10").Value = "M/D/YYYY h:mm:ss[.fff] TT"
$ws.Range("C11").Value = "h:mm:ss TT"

# Status fill changes from red (not started) to yellow (in progress)
$ws.Range("B8:C11").Interior.Color = 65535

# Keep the wrap-text formatting on C11
$ws.Range("C11").WrapText = $true

# Move the active selection to the last cell touched
$ws.Range("C11").Select()
